$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.734.96"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "3.340.70"
$ws.Range("E3").Value = "  -4.03%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.16"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.08"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +4.16%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "3.337.31"

$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "3.924.97"
$ws.Range("E13").Value = "  -3.89%  "

$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.71"
$ws.Range("E15").Value = "  -4.69%  "

$ws.Range("D16").Value = "65.753.13"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000169"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").Value = "3.322.07"
$ws.Range("E18").Value = "  -4.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.73"
$ws.Range("E19").Value = "  -3.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.44"
$ws.Range("E20").Value = "  -3.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.86"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.44"
$ws.Range("E22").Value = "  -4.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.76"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.67"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.81"
$ws.Range("E33").Value = "  -5.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.88"
$ws.Range("E34").Value = "  -3.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.22"
$ws.Range("E35").Value = "  -5.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.61"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.849"
$ws.Range("E38").Value = "  -4.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.19"
$ws.Range("E39").Value = "  -7.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -2.44%  "

$ws.Range("D42").Value = "2.664.83"
$ws.Range("E42").Value = "  -5.85%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.28"
$ws.Range("E43").Value = "  -3.92%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.21"
$ws.Range("E44").Value = "  -3.62%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "335.86"
$ws.Range("E45").Value = "  +8.05%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0666"
$ws.Range("E46").Value = "  -2.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.75"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.41"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0279"
$ws.Range("E49").Value = "  -3.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  +2.50%  "

$ws.Range("E51").Value = "  -1.04%  "
